# Re-sync the OH "top remaining" scraper sheet with the latest scrape pass.
#
# Two kinds of change landed together:
#   1) Every row that was last scraped on 2019-03-07 gets bumped to the new
#      scrape date, 2019-03-12 (rows still dated 2019-02-13 are untouched --
#      those games weren't re-scraped this pass).
#   2) A handful of rows shuffled position in the scraper's output between
#      runs, so their GAME NAME / GAME NUMBER / TOP PRIZES REMAINING values
#      need to be realigned to match the new scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "2019-03-07"
$newDate = "2019-03-12"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dateCell = $ws.Cells.Item($r, 6)
    if ($dateCell.Value2 -eq $oldDate) {
        # Leading apostrophe forces text so Excel doesn't re-type the cell
        # as a date (the column stores scrape dates as plain text).
        $dateCell.Formula = "'" + $newDate
    }
}

# Rows whose GAME NAME / GAME NUMBER / TOP PRIZES REMAINING were realigned
# to the freshly scraped order.
$rowFixups = @(
    @{ Row = 2;  Name = "Double Doubler";              Number = 490; Prizes = 14  },
    @{ Row = 3;  Name = "Luck of the Irish Tripler";   Number = 491; Prizes = 18  },
    @{ Row = 5;                                                      Prizes = 25  },
    @{ Row = 9;  Name = "Snow Me the Money";           Number = 485; Prizes = 12  },
    @{ Row = 10; Name = "Holiday Cash";                Number = 480; Prizes = 116 },
    @{ Row = 18; Name = "Fireball Bingo";              Number = 502; Prizes = 9   },
    @{ Row = 19; Name = "Decade of Dollars";           Number = 497; Prizes = 5   },
    @{ Row = 21; Name = "Power Play Cashword";         Number = 462; Prizes = 2   },
    @{ Row = 22; Name = "Blazing Hot Cash";            Number = 457; Prizes = 5   },
    @{ Row = 23; Name = "I Love Cash";                 Number = 492; Prizes = 6   },
    @{ Row = 25;                                                     Prizes = 43  },
    @{ Row = 34; Name = "Skee-Ball";                   Number = 474; Prizes = 1   },
    @{ Row = 35; Name = "Bingo Plus";                  Number = 404               },
    @{ Row = 36; Name = "Cash Wheel";                  Number = 498; Prizes = 2   },
    @{ Row = 61;                                                     Prizes = 2   }
)

foreach ($fix in $rowFixups) {
    $r = $fix.Row
    if ($fix.ContainsKey("Name"))   { $ws.Cells.Item($r, 3).Formula = $fix.Name }
    if ($fix.ContainsKey("Number")) { $ws.Cells.Item($r, 4).Formula = $fix.Number }
    if ($fix.ContainsKey("Prizes")) { $ws.Cells.Item($r, 5).Formula = $fix.Prizes }
}
